# Append a new "Installation bootstrap" section at the very end of the
# document, just before the final (empty) trailing paragraph.
#
# The new content consists of:
#   - one empty paragraph (spacer)
#   - "Installation bootstrap"
#   - "Verifier la version :"
#   - "npm view bootstrap version"
#   - "installer bootstrap en lignes de commandes :"
#   - "npm install bootstrap"
#
# We build the exact OOXML for these paragraphs (including the
# w:proofErr spell/grammar-check markers that Word itself would emit)
# and inject it with Range.InsertXML so the resulting markup matches
# precisely, instead of relying on autocorrect/spellcheck side effects
# that this headless host does not simulate.

$d = $word.ActiveDocument

# Locate the last paragraph in the document (the trailing empty
# paragraph that sits right before the sectPr) and collapse a range to
# its very start - that's exactly where the new paragraphs must land.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(1)

$newParagraphsXml = '<w:p/>' +
  '<w:p><w:r><w:t xml:space="preserve">Installation </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>bootstrap</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Verifier</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> la version :</w:t></w:r></w:p>' +
  '<w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>npm</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>view</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bootstrap</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> version</w:t></w:r></w:p>' +
  '<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>installer</w:t></w:r><w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bootstrap</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en lignes de commandes</w:t></w:r>' +
  '<w:r><w:t> :</w:t></w:r></w:p>' +
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>install</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> bootstrap</w:t></w:r></w:p>'

$flatOpcXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' + $newParagraphsXml + '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($flatOpcXml)
